$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.67774772644043
$ws.Range("B1").Value = 1.939854264259338
$ws.Range("C1").Value = 2.024612903594971
$ws.Range("D1").Value = 2.515945911407471
$ws.Range("E1").Value = 3.562259197235107
